# Insert a new blank row above row 142 (pushing rows 142+ down by one),
# copying the cell formatting from row 141 so the new blank row matches
# the style of the row immediately above it (as Excel does on row insert).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(142).Insert(-4121, 0) | Out-Null
$ws.Range("A141:D141").Copy() | Out-Null
$ws.Range("A142:D142").PasteSpecial(-4122) | Out-Null

# Restore the selection / active cell as recorded after the edit.
$ws.Range("A145").Select() | Out-Null
